$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 63.81813833333333
$ws.Range("H2").Value = 191.454415
$ws.Range("I2").Value = 0.5585681932726833
$ws.Range("J2").Value = 0.5585681932726834
$ws.Range("M2").Value = 22.495411
$ws.Range("N2").Value = 67.486233
$ws.Range("O2").Value = 0.08292345339295874
$ws.Range("P2").Value = 0.08292345339295874
$ws.Range("Q2").Value = 1435.615251063188
$ws.Range("R2").Value = 12920.53725956869
$ws.Range("S2").Value = 0.04631840354163652
$ws.Range("T2").Value = 0.04631840354163653

# Row 3
$ws.Range("G3").Value = 63.81813833333333
$ws.Range("H3").Value = 191.454415
$ws.Range("I3").Value = 0.5585681932726833
$ws.Range("J3").Value = 0.5585681932726834
$ws.Range("M3").Value = 82.64333833333332
$ws.Range("O3").Value = 0.3046430676248896
$ws.Range("P3").Value = 0.3046430676248896
$ws.Range("Q3").Value = 5274.143998085135
$ws.Range("R3").Value = 47467.29598276621
$ws.Range("S3").Value = 0.1701639278762825
$ws.Range("T3").Value = 0.1701639278762825

# Row 4
$ws.Range("G4").Value = 63.81813833333333
$ws.Range("H4").Value = 191.454415
$ws.Range("I4").Value = 0.5585681932726833
$ws.Range("J4").Value = 0.5585681932726834
$ws.Range("M4").Value = 79.32606499999999
$ws.Range("N4").Value = 237.978195
$ws.Range("O4").Value = 0.2924148064631633
$ws.Range("P4").Value = 0.2924148064631633
$ws.Range("Q4").Value = 5062.441789608991
$ws.Range("R4").Value = 45561.97610648091
$ws.Range("S4").Value = 0.1633336101323105
$ws.Range("T4").Value = 0.1633336101323105

# Row 5
$ws.Range("G5").Value = 63.81813833333333
$ws.Range("H5").Value = 191.454415
$ws.Range("I5").Value = 0.5585681932726833
$ws.Range("J5").Value = 0.5585681932726834
$ws.Range("M5").Value = 14.467164
$ws.Range("N5").Value = 43.401492
$ws.Range("O5").Value = 0.05332941903938943
$ws.Range("P5").Value = 0.05332941903938943
$ws.Range("Q5").Value = 923.2674734430199
$ws.Range("R5").Value = 8309.407260987178
$ws.Range("S5").Value = 0.02978811724111359
$ws.Range("T5").Value = 0.0297881172411136

# Row 6
$ws.Range("G6").Value = 63.81813833333333
$ws.Range("H6").Value = 191.454415
$ws.Range("I6").Value = 0.5585681932726833
$ws.Range("J6").Value = 0.5585681932726834
$ws.Range("M6").Value = 72.34725666666667
$ws.Range("N6").Value = 217.04177
$ws.Range("O6").Value = 0.2666892534795989
$ws.Range("P6").Value = 0.2666892534795989
$ws.Range("Q6").Value = 4617.067233990505
$ws.Range("R6").Value = 41553.60510591455
$ws.Range("S6").Value = 0.1489641344813402
$ws.Range("T6").Value = 0.1489641344813403

# Row 7
$ws.Range("H7").Value = 44.084775
$ws.Range("I7").Value = 0.1286173166734377
$ws.Range("J7").Value = 0.1286173166734377
$ws.Range("M7").Value = 22.495411
$ws.Range("N7").Value = 67.486233
$ws.Range("O7").Value = 0.08292345339295874
$ws.Range("P7").Value = 0.08292345339295874
$ws.Range("Q7").Value = 330.568377489175
$ws.Range("R7").Value = 2975.115397402575
$ws.Range("S7").Value = 0.01066539206469722
$ws.Range("T7").Value = 0.01066539206469723

# Row 8
$ws.Range("H8").Value = 44.084775
$ws.Range("I8").Value = 0.1286173166734377
$ws.Range("J8").Value = 0.1286173166734377
$ws.Range("M8").Value = 82.64333833333332
$ws.Range("O8").Value = 0.3046430676248896
$ws.Range("P8").Value = 0.3046430676248896
$ws.Range("Q8").Value = 1214.437658557958
$ws.Range("R8").Value = 10929.93892702162
$ws.Range("S8").Value = 0.03918237390107792
$ws.Range("T8").Value = 0.03918237390107793

# Row 9
$ws.Range("H9").Value = 44.084775
$ws.Range("I9").Value = 0.1286173166734377
$ws.Range("J9").Value = 0.1286173166734377
$ws.Range("M9").Value = 79.32606499999999
$ws.Range("N9").Value = 237.978195
$ws.Range("O9").Value = 0.2924148064631633
$ws.Range("P9").Value = 0.2924148064631633
$ws.Range("Q9").Value = 1165.690575720125
$ws.Range("R9").Value = 10491.21518148112
$ws.Range("S9").Value = 0.03760960776287466
$ws.Range("T9").Value = 0.03760960776287468

# Row 10
$ws.Range("H10").Value = 44.084775
$ws.Range("I10").Value = 0.1286173166734377
$ws.Range("J10").Value = 0.1286173166734377
$ws.Range("M10").Value = 14.467164
$ws.Range("N10").Value = 43.401492
$ws.Range("O10").Value = 0.05332941903938943
$ws.Range("P10").Value = 0.05332941903938943
$ws.Range("Q10").Value = 212.5938899427
$ws.Range("R10").Value = 1913.3450094843
$ws.Range("S10").Value = 0.006859086776599607
$ws.Range("T10").Value = 0.006859086776599609

# Row 11
$ws.Range("H11").Value = 44.084775
$ws.Range("I11").Value = 0.1286173166734377
$ws.Range("J11").Value = 0.1286173166734377
$ws.Range("M11").Value = 72.34725666666667
$ws.Range("N11").Value = 217.04177
$ws.Range("O11").Value = 0.2666892534795989
$ws.Range("P11").Value = 0.2666892534795989
$ws.Range("Q11").Value = 1063.137510672417
$ws.Range("R11").Value = 9568.237596051751
$ws.Range("S11").Value = 0.03430085616818826
$ws.Range("T11").Value = 0.03430085616818827

# Row 12
$ws.Range("G12").Value = 17.02115633333333
$ws.Range("H12").Value = 51.063469
$ws.Range("I12").Value = 0.1489776541406249
$ws.Range("J12").Value = 0.1489776541406249
$ws.Range("M12").Value = 22.495411
$ws.Range("N12").Value = 67.486233
$ws.Range("O12").Value = 0.08292345339295874
$ws.Range("P12").Value = 0.08292345339295874
$ws.Range("Q12").Value = 382.8979074135864
$ws.Range("R12").Value = 3446.081166722277
$ws.Range("S12").Value = 0.01235374155972244
$ws.Range("T12").Value = 0.01235374155972244

# Row 13
$ws.Range("G13").Value = 17.02115633333333
$ws.Range("H13").Value = 51.063469
$ws.Range("I13").Value = 0.1489776541406249
$ws.Range("J13").Value = 0.1489776541406249
$ws.Range("M13").Value = 82.64333833333332
$ws.Range("O13").Value = 0.3046430676248896
$ws.Range("P13").Value = 0.3046430676248896
$ws.Range("Q13").Value = 1406.685181680226
$ws.Range("R13").Value = 12660.16663512203
$ws.Range("S13").Value = 0.04538500956495981
$ws.Range("T13").Value = 0.04538500956495983

# Row 14
$ws.Range("G14").Value = 17.02115633333333
$ws.Range("H14").Value = 51.063469
$ws.Range("I14").Value = 0.1489776541406249
$ws.Range("J14").Value = 0.1489776541406249
$ws.Range("M14").Value = 79.32606499999999
$ws.Range("N14").Value = 237.978195
$ws.Range("O14").Value = 0.2924148064631633
$ws.Range("P14").Value = 0.2924148064631633
$ws.Range("Q14").Value = 1350.221353673162
$ws.Range("R14").Value = 12151.99218305846
$ws.Range("S14").Value = 0.04356327190286691
$ws.Range("T14").Value = 0.04356327190286692

# Row 15
$ws.Range("G15").Value = 17.02115633333333
$ws.Range("H15").Value = 51.063469
$ws.Range("I15").Value = 0.1489776541406249
$ws.Range("J15").Value = 0.1489776541406249
$ws.Range("M15").Value = 14.467164
$ws.Range("N15").Value = 43.401492
$ws.Range("O15").Value = 0.05332941903938943
$ws.Range("P15").Value = 0.05332941903938943
$ws.Range("Q15").Value = 246.247860143972
$ws.Range("R15").Value = 2216.230741295748
$ws.Range("S15").Value = 0.007944891745170617
$ws.Range("T15").Value = 0.007944891745170617

# Row 16
$ws.Range("G16").Value = 17.02115633333333
$ws.Range("H16").Value = 51.063469
$ws.Range("I16").Value = 0.1489776541406249
$ws.Range("J16").Value = 0.1489776541406249
$ws.Range("M16").Value = 72.34725666666667
$ws.Range("N16").Value = 217.04177
$ws.Range("O16").Value = 0.2666892534795989
$ws.Range("P16").Value = 0.2666892534795989
$ws.Range("Q16").Value = 1231.433966011126
$ws.Range("R16").Value = 11082.90569410013
$ws.Range("S16").Value = 0.03973073936790513
$ws.Range("T16").Value = 0.03973073936790514

# Row 17
$ws.Range("G17").Value = 0.7288956666666667
$ws.Range("H17").Value = 2.186687
$ws.Range("I17").Value = 0.006379658608775693
$ws.Range("J17").Value = 0.006379658608775693
$ws.Range("M17").Value = 22.495411
$ws.Range("N17").Value = 67.486233
$ws.Range("O17").Value = 0.08292345339295874
$ws.Range("P17").Value = 0.08292345339295874
$ws.Range("Q17").Value = 16.39680759778567
$ws.Range("R17").Value = 147.571268380071
$ws.Range("S17").Value = 0.0005290233233077991
$ws.Range("T17").Value = 0.0005290233233077991

# Row 18
$ws.Range("G18").Value = 0.7288956666666667
$ws.Range("H18").Value = 2.186687
$ws.Range("I18").Value = 0.006379658608775693
$ws.Range("J18").Value = 0.006379658608775693
$ws.Range("M18").Value = 82.64333833333332
$ws.Range("O18").Value = 0.3046430676248896
$ws.Range("P18").Value = 0.3046430676248896
$ws.Range("Q18").Value = 60.23837119003388
$ws.Range("R18").Value = 542.1453407103049
$ws.Range("S18").Value = 0.001943518768976963
$ws.Range("T18").Value = 0.001943518768976963

# Row 19
$ws.Range("G19").Value = 0.7288956666666667
$ws.Range("H19").Value = 2.186687
$ws.Range("I19").Value = 0.006379658608775693
$ws.Range("J19").Value = 0.006379658608775693
$ws.Range("M19").Value = 79.32606499999999
$ws.Range("N19").Value = 237.978195
$ws.Range("O19").Value = 0.2924148064631633
$ws.Range("P19").Value = 0.2924148064631633
$ws.Range("Q19").Value = 57.82042503221832
$ws.Range("R19").Value = 520.383825289965
$ws.Range("S19").Value = 0.001865506637386198
$ws.Range("T19").Value = 0.001865506637386198

# Row 20
$ws.Range("G20").Value = 0.7288956666666667
$ws.Range("H20").Value = 2.186687
$ws.Range("I20").Value = 0.006379658608775693
$ws.Range("J20").Value = 0.006379658608775693
$ws.Range("M20").Value = 14.467164
$ws.Range("N20").Value = 43.401492
$ws.Range("O20").Value = 0.05332941903938943
$ws.Range("P20").Value = 0.05332941903938943
$ws.Range("Q20").Value = 10.545053148556
$ws.Range("R20").Value = 94.905478337004
$ws.Range("S20").Value = 0.0003402234872756471
$ws.Range("T20").Value = 0.0003402234872756471

# Row 21
$ws.Range("G21").Value = 0.7288956666666667
$ws.Range("H21").Value = 2.186687
$ws.Range("I21").Value = 0.006379658608775693
$ws.Range("J21").Value = 0.006379658608775693
$ws.Range("M21").Value = 72.34725666666667
$ws.Range("N21").Value = 217.04177
$ws.Range("O21").Value = 0.2666892534795989
$ws.Range("P21").Value = 0.2666892534795989
$ws.Range("Q21").Value = 52.73360187955445
$ws.Range("R21").Value = 474.60241691599
$ws.Range("S21").Value = 0.001701386391829086
$ws.Range("T21").Value = 0.001701386391829086

# Row 22
$ws.Range("G22").Value = 17.989968
$ws.Range("H22").Value = 53.969904
$ws.Range("I22").Value = 0.1574571773044783
$ws.Range("J22").Value = 0.1574571773044783
$ws.Range("M22").Value = 22.495411
$ws.Range("N22").Value = 67.486233
$ws.Range("O22").Value = 0.08292345339295874
$ws.Range("P22").Value = 0.08292345339295874
$ws.Range("Q22").Value = 404.691724036848
$ws.Range("R22").Value = 3642.225516331632
$ws.Range("S22").Value = 0.01305689290359475
$ws.Range("T22").Value = 0.01305689290359475

# Row 23
$ws.Range("G23").Value = 17.989968
$ws.Range("H23").Value = 53.969904
$ws.Range("I23").Value = 0.1574571773044783
$ws.Range("J23").Value = 0.1574571773044783
$ws.Range("M23").Value = 82.64333833333332
$ws.Range("O23").Value = 0.3046430676248896
$ws.Range("P23").Value = 0.3046430676248896
$ws.Range("Q23").Value = 1486.75101202984
$ws.Range("R23").Value = 13380.75910826856
$ws.Range("S23").Value = 0.04796823751359242
$ws.Range("T23").Value = 0.04796823751359242

# Row 24
$ws.Range("G24").Value = 17.989968
$ws.Range("H24").Value = 53.969904
$ws.Range("I24").Value = 0.1574571773044783
$ws.Range("J24").Value = 0.1574571773044783
$ws.Range("M24").Value = 79.32606499999999
$ws.Range("N24").Value = 237.978195
$ws.Range("O24").Value = 0.2924148064631633
$ws.Range("P24").Value = 0.2924148064631633
$ws.Range("Q24").Value = 1427.07337091592
$ws.Range("R24").Value = 12843.66033824328
$ws.Range("S24").Value = 0.04604281002772501
$ws.Range("T24").Value = 0.04604281002772501

# Row 25
$ws.Range("G25").Value = 17.989968
$ws.Range("H25").Value = 53.969904
$ws.Range("I25").Value = 0.1574571773044783
$ws.Range("J25").Value = 0.1574571773044783
$ws.Range("M25").Value = 14.467164
$ws.Range("N25").Value = 43.401492
$ws.Range("O25").Value = 0.05332941903938943
$ws.Range("P25").Value = 0.05332941903938943
$ws.Range("Q25").Value = 260.263817410752
$ws.Range("R25").Value = 2342.374356696768
$ws.Range("S25").Value = 0.008397099789229961
$ws.Range("T25").Value = 0.008397099789229961

# Row 26
$ws.Range("G26").Value = 17.989968
$ws.Range("H26").Value = 53.969904
$ws.Range("I26").Value = 0.1574571773044783
$ws.Range("J26").Value = 0.1574571773044783
$ws.Range("M26").Value = 72.34725666666667
$ws.Range("N26").Value = 217.04177
$ws.Range("O26").Value = 0.2666892534795989
$ws.Range("P26").Value = 0.2666892534795989
$ws.Range("Q26").Value = 1301.52483232112
$ws.Range("R26").Value = 11713.72349089008
$ws.Range("S26").Value = 0.04199213707033615
$ws.Range("T26").Value = 0.04199213707033616
